$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 438
$ws.Range("I19").Value = 478.66666
$ws.Range("J19").Value = 403.14285
$ws.Range("K19").Value = 478.66666
$ws.Range("L19").Value = 403.14285
$ws.Range("M19").Value = -303.66666
$ws.Range("N19").Value = -753.14285

$ws.Range("H43").Value = 5567255.5
$ws.Range("I43").Value = 33800.332
$ws.Range("J43").Value = 7938736.5
$ws.Range("K43").Value = 33800.332
$ws.Range("L43").Value = 7938736.5
$ws.Range("M43").Value = -33731.332
$ws.Range("N43").Value = -7938874.5

$ws.Range("H62").Value = 11114960
$ws.Range("I62").Value = 13892699
$ws.Range("J62").Value = 4003
$ws.Range("K62").Value = 13892699
$ws.Range("L62").Value = 4003
$ws.Range("M62").Value = -13892075
$ws.Range("N62").Value = -5251

$ws.Range("H65").Value = 11114960
$ws.Range("I65").Value = 13892699
$ws.Range("J65").Value = 4003
$ws.Range("K65").Value = 69463495
$ws.Range("L65").Value = 20015
$ws.Range("M65").Value = -69460375
$ws.Range("N65").Value = -26255

$ws.Range("H97").Value = 397
$ws.Range("J97").Value = 379.83334
$ws.Range("L97").Value = 1139.50002
$ws.Range("N97").Value = -2131.50002

$ws.Range("H98").Value = 7047.4614
$ws.Range("I98").Value = 7628.909
$ws.Range("K98").Value = 7628.909
$ws.Range("M98").Value = -6130.909

$ws.Range("H121").Value = 1747.5
$ws.Range("J121").Value = 1830
$ws.Range("L121").Value = 5490
$ws.Range("N121").Value = -8984

$ws.Range("H122").Value = 7047.4614
$ws.Range("I122").Value = 7628.909
$ws.Range("K122").Value = 22886.727
$ws.Range("M122").Value = -20436.727

$ws.Range("H129").Value = 804.19354
$ws.Range("I129").Value = 552.4286
$ws.Range("J129").Value = 877.625
$ws.Range("K129").Value = 1657.2858
$ws.Range("L129").Value = 2632.875
$ws.Range("M129").Value = 3342.7142
$ws.Range("N129").Value = -12632.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 147.27272
$ws.Range("I5").Value = 120.14286
$ws.Range("J5").Value = 194.75
$ws.Range("K5").Value = 120.14286
$ws.Range("L5").Value = 194.75
$ws.Range("M5").Value = -8.142859999999999
$ws.Range("N5").Value = -418.75

$ws.Range("H32").Value = 7339.579
$ws.Range("I32").Value = 7267.7026
$ws.Range("K32").Value = 7267.7026
$ws.Range("M32").Value = -6980.7026

$ws.Range("H74").Value = 1006.3333
$ws.Range("I74").Value = 1022.2
$ws.Range("J74").Value = 927
$ws.Range("K74").Value = 1022.2
$ws.Range("L74").Value = 927
$ws.Range("M74").Value = -148.2
$ws.Range("N74").Value = -2675

$ws.Range("H77").Value = 1006.3333
$ws.Range("I77").Value = 1022.2
$ws.Range("J77").Value = 927
$ws.Range("K77").Value = 5111
$ws.Range("L77").Value = 4635
$ws.Range("M77").Value = -743
$ws.Range("N77").Value = -13371

$ws.Range("H132").Value = 4046.5293
$ws.Range("I132").Value = 4580.125
$ws.Range("K132").Value = 13740.375
$ws.Range("M132").Value = -11210.375

$ws.Range("H135").Value = 62714.5
$ws.Range("J135").Value = 62714.5
$ws.Range("L135").Value = 62714.5
$ws.Range("N135").Value = -72854.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 147.27272
$ws.Range("I4").Value = 120.14286
$ws.Range("J4").Value = 194.75
$ws.Range("K4").Value = 120.14286
$ws.Range("L4").Value = 194.75
$ws.Range("M4").Value = -5.142859999999999
$ws.Range("N4").Value = -424.75

$ws.Range("H64").Value = 460.66666
$ws.Range("I64").Value = 329.83334
$ws.Range("J64").Value = 591.5
$ws.Range("K64").Value = 329.83334
$ws.Range("L64").Value = 591.5
$ws.Range("M64").Value = -104.83334
$ws.Range("N64").Value = -1041.5

$ws.Range("H67").Value = 460.66666
$ws.Range("I67").Value = 329.83334
$ws.Range("J67").Value = 591.5
$ws.Range("K67").Value = 329.83334
$ws.Range("L67").Value = 591.5
$ws.Range("M67").Value = 450.16666
$ws.Range("N67").Value = -2151.5

$ws.Range("H94").Value = 14706621
$ws.Range("I94").Value = 20834016
$ws.Range("J94").Value = 872
$ws.Range("K94").Value = 20834016
$ws.Range("L94").Value = 872
$ws.Range("M94").Value = -20833565
$ws.Range("N94").Value = -1774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1645.1936
$ws.Range("I31").Value = 1864.4667
$ws.Range("K31").Value = 1864.4667
$ws.Range("M31").Value = -1569.4667

$ws.Range("H34").Value = 1645.1936
$ws.Range("I34").Value = 1864.4667
$ws.Range("K34").Value = 1864.4667
$ws.Range("M34").Value = -1662.4667

$ws.Range("H58").Value = 1539.7084
$ws.Range("I58").Value = 1208.2
$ws.Range("K58").Value = 1208.2
$ws.Range("M58").Value = -1005.2

$ws.Range("H105").Value = 728
$ws.Range("I105").Value = 723.2
$ws.Range("K105").Value = 723.2
$ws.Range("M105").Value = 1023.8

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H122").Value = 4266.857
$ws.Range("I122").Value = 4503.154
$ws.Range("K122").Value = 13509.462
$ws.Range("M122").Value = -11059.462

$ws.Range("H132").Value = 2377.3333
$ws.Range("I132").Value = 2133.5
$ws.Range("J132").Value = 3157.6
$ws.Range("K132").Value = 6400.5
$ws.Range("L132").Value = 9472.799999999999
$ws.Range("M132").Value = -3870.5
$ws.Range("N132").Value = -14532.8

$ws.Range("H136").Value = 1539.7084
$ws.Range("I136").Value = 1208.2
$ws.Range("K136").Value = 3624.6
$ws.Range("M136").Value = -1074.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2234.2144
$ws.Range("J34").Value = 2564.9167
$ws.Range("L34").Value = 7694.750100000001
$ws.Range("N34").Value = -7862.750100000001

$ws.Range("H97").Value = 1100
$ws.Range("J97").Value = 1500
$ws.Range("L97").Value = 4500
$ws.Range("N97").Value = -5492

$ws.Range("H131").Value = 26317422
$ws.Range("J131").Value = 1830.5758
$ws.Range("L131").Value = 5491.7274
$ws.Range("N131").Value = -15571.7274

$ws.Range("H132").Value = 879.26086
$ws.Range("I132").Value = 891.5333000000001
$ws.Range("J132").Value = 856.25
$ws.Range("K132").Value = 8023.7997
$ws.Range("L132").Value = 7706.25
$ws.Range("M132").Value = -5493.7997
$ws.Range("N132").Value = -12766.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 225000000
$ws.Range("I70").Value = 250000000
$ws.Range("K70").Value = 250000000
$ws.Range("M70").Value = -249999730

$ws.Range("H73").Value = 225000000
$ws.Range("I73").Value = 250000000
$ws.Range("K73").Value = 250000000
$ws.Range("M73").Value = -249999064

$ws.Range("H102").Value = 827.41174
$ws.Range("I102").Value = 718.7273
$ws.Range("J102").Value = 1026.6666
$ws.Range("K102").Value = 718.7273
$ws.Range("L102").Value = 1026.6666
$ws.Range("M102").Value = 903.2727
$ws.Range("N102").Value = -4270.6666

$ws.Range("H122").Value = 2266.5806
$ws.Range("I122").Value = 2348.4783
$ws.Range("K122").Value = 7045.4349
$ws.Range("M122").Value = -4595.4349

$ws.Range("H126").Value = 1975.2142
$ws.Range("I126").Value = 1754.875
$ws.Range("J126").Value = 2269
$ws.Range("K126").Value = 5264.625
$ws.Range("L126").Value = 6807
$ws.Range("M126").Value = -2794.625
$ws.Range("N126").Value = -11747

$ws.Range("H132").Value = 4555.385
$ws.Range("I132").Value = 5951.5
$ws.Range("J132").Value = 3358.7144
$ws.Range("K132").Value = 17854.5
$ws.Range("L132").Value = 10076.1432
$ws.Range("M132").Value = -15324.5
$ws.Range("N132").Value = -15136.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 10857.143
$ws.Range("I56").Value = 8000
$ws.Range("K56").Value = 8000
$ws.Range("M56").Value = -7309

$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352

$ws.Range("H122").Value = 50001840
$ws.Range("I122").Value = 83334800
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 250004400
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -250001950
$ws.Range("N122").Value = -12100

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 7828.5
$ws.Range("I61").Value = 4600
$ws.Range("J61").Value = 11057
$ws.Range("K61").Value = 4600
$ws.Range("L61").Value = 11057
$ws.Range("M61").Value = -4308
$ws.Range("N61").Value = -11641

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H126").Value = 38462268
$ws.Range("I126").Value = 43478890
$ws.Range("J126").Value = 1518.3334
$ws.Range("K126").Value = 130436670
$ws.Range("L126").Value = 4555.0002
$ws.Range("M126").Value = -130434200
$ws.Range("N126").Value = -9495.0002

$ws.Range("H132").Value = 1711.3
$ws.Range("I132").Value = 1278.6923
$ws.Range("J132").Value = 2514.7144
$ws.Range("K132").Value = 3836.0769
$ws.Range("L132").Value = 7544.1432
$ws.Range("M132").Value = -1306.0769
$ws.Range("N132").Value = -12604.1432
